$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1) Sheet "研究生获得奖励" (graduate-student awards):
#    Insert a new row for Yin Ting's CCUS paper award (dated 2024)
#    before the current row 6 ("研究生国家奖学金"), which pushes the
#    following two rows down by one.
# -------------------------------------------------------------------
$wsAwards = $wb.Worksheets.Item("研究生获得奖励")

$wsAwards.Rows.Item(6).Insert()

$wsAwards.Range("A6").Value = 5
$wsAwards.Range("B6").Value = "尹亭"
$wsAwards.Range("C6").Value = "全国碳捕集、利用与封存（CCUS）技术研讨会优秀论文二等奖"
$wsAwards.Range("D6").Value = 2024
$wsAwards.Rows.Item(6).RowHeight = 16.15

# Renumber the "编号" column for the rows that shifted down.
$wsAwards.Range("A7").Value = 6
$wsAwards.Range("A8").Value = 7
$wsAwards.Range("A9").Value = 8

# Update the view's current selection for this sheet.
$wsAwards.Range("F12").Select()

# -------------------------------------------------------------------
# 2) Sheet "发表论文" (published papers): append Yin et al. (2024)
#    Energy paper as a new row with its hyperlink.
# -------------------------------------------------------------------
$wsPapers = $wb.Worksheets.Item("发表论文")

$wsPapers.Range("A40").Value = "Yin, T., Chen, S., Wang, G., Tan, Y., Teng, F., Zhang, Q., 2024. Can Subsidy Policies Achieve Fuel Cell Logistics Vehicle (FCLV) Promotion Targets? Evidence from the Beijing-Tianjin-Hebei Fuel Cell Vehicle Demonstration City Cluster in China. Energy 133270. https://doi.org/10.1016/j.energy.2024.133270"

$wsPapers.Hyperlinks.Add($wsPapers.Range("A40"), "https://doi.org/10.1016/j.energy.2024.133270", "", "", "https://doi.org/10.1016/j.energy.2024.133270")

# Re-apply the same "hyperlink" cell formatting used by the other rows
# in this column (keeps the cell style consistent with the rest of the
# sheet instead of a freshly-synthesised duplicate style).
$wsPapers.Range("A2").Copy()
$wsPapers.Range("A40").PasteSpecial(-4122)

# -------------------------------------------------------------------
# 3) Make "发表论文" the active sheet/tab again, with its view scrolled
#    and the selection parked on D44.
# -------------------------------------------------------------------
$wsPapers.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 25
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$wsPapers.Range("D44").Select()
